$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "'2026-02-05 17:18:06"
$ws.Range("O2").Value = "'-1.1 °C"
$ws.Range("E3").Value = "'2026-02-05 17:18:08"
$ws.Range("E4").Value = "'2026-02-05 17:18:11"
$ws.Range("J4").Value = "'990.6 hPa"
$ws.Range("L4").Value = "'81.7 km/h - 283º 16:31 TU"
$ws.Range("O4").Value = "'10.8 °C"
$ws.Range("E5").Value = "'2026-02-05 17:18:14"
$ws.Range("H5").Value = "'74%"
$ws.Range("J5").Value = "'990.6 hPa"
$ws.Range("O5").Value = "'9.2 °C"
$ws.Range("E6").Value = "'2026-02-05 17:18:17"
$ws.Range("J6").Value = "'992.2 hPa"
$ws.Range("K6").Value = "'3.8 MJ/m2"
$ws.Range("O6").Value = "'12.5 °C"
$ws.Range("E7").Value = "'2026-02-05 17:18:19"
$ws.Range("J7").Value = "'992.0 hPa"
$ws.Range("E8").Value = "'2026-02-05 17:18:22"
$ws.Range("K8").Value = "'5.6 MJ/m2"
$ws.Range("O8").Value = "'8.0 °C"
$ws.Range("E9").Value = "'2026-02-05 17:18:25"
$ws.Range("O9").Value = "'1.9 °C"
$ws.Range("E10").Value = "'2026-02-05 17:18:28"
$ws.Range("O10").Value = "'7.1 °C"
$ws.Range("E11").Value = "'2026-02-05 17:18:30"
$ws.Range("J11").Value = "'995.5 hPa"
$ws.Range("L11").Value = "'20.5 km/h - 204º 16:42 TU"
$ws.Range("O11").Value = "'0.1 °C"
$ws.Range("E12").Value = "'2026-02-05 17:18:33"
$ws.Range("M12").Value = "'15.8 °C 16:43 TU"
$ws.Range("O12").Value = "'9.2 °C"
$ws.Range("E13").Value = "'2026-02-05 17:18:36"
$ws.Range("O13").Value = "'7.5 °C"
$ws.Range("E14").Value = "'2026-02-05 17:18:39"
$ws.Range("I14").Value = "'5.6 mm"
$ws.Range("E15").Value = "'2026-02-05 17:18:41"
$ws.Range("H15").Value = "'86%"
$ws.Range("J15").Value = "'991.1 hPa"
$ws.Range("K15").Value = "'6.7 MJ/m2"
$ws.Range("M15").Value = "'16.5 °C 16:58 TU"
$ws.Range("O15").Value = "'7.2 °C"
$ws.Range("E16").Value = "'2026-02-05 17:18:44"
$ws.Range("M16").Value = "'6.3 °C 16:34 TU"
$ws.Range("E17").Value = "'2026-02-05 17:18:47"
$ws.Range("I17").Value = "'8.1 mm"
$ws.Range("J17").Value = "'995.5 hPa"
$ws.Range("E18").Value = "'2026-02-05 17:18:50"
$ws.Range("O18").Value = "'-4.4 °C"
$ws.Range("E19").Value = "'2026-02-05 17:18:53"
$ws.Range("I19").Value = "'7.7 mm"
$ws.Range("J19").Value = "'992.6 hPa"
$ws.Range("O19").Value = "'7.2 °C"
$ws.Range("E20").Value = "'2026-02-05 17:18:55"
$ws.Range("H20").Value = "'73%"
$ws.Range("K20").Value = "'1.3 MJ/m2"
$ws.Range("L20").Value = "'46.8 km/h - 251º 16:58 TU"
$ws.Range("O20").Value = "'-1.6 °C"
$ws.Range("E21").Value = "'2026-02-05 17:18:58"
$ws.Range("J21").Value = "'991.4 hPa"
$ws.Range("O21").Value = "'5.6 °C"
$ws.Range("E22").Value = "'2026-02-05 17:19:01"
$ws.Range("H22").Value = "'91%"
$ws.Range("M22").Value = "'14.7 °C 16:40 TU"
$ws.Range("O22").Value = "'7.8 °C"
$ws.Range("E23").Value = "'2026-02-05 17:19:04"
$ws.Range("J23").Value = "'990.5 hPa"
$ws.Range("K23").Value = "'3.5 MJ/m2"
$ws.Range("O23").Value = "'7.9 °C"
$ws.Range("E24").Value = "'2026-02-05 17:19:06"
$ws.Range("J24").Value = "'989.5 hPa"
$ws.Range("O24").Value = "'10.1 °C"
$ws.Range("E25").Value = "'2026-02-05 17:19:09"
$ws.Range("J25").Value = "'994.6 hPa"
$ws.Range("E26").Value = "'2026-02-05 17:19:11"
$ws.Range("H26").Value = "'75%"
$ws.Range("O26").Value = "'-1.0 °C"
$ws.Range("E27").Value = "'2026-02-05 17:19:14"
$ws.Range("J27").Value = "'990.9 hPa"
$ws.Range("O27").Value = "'8.2 °C"
$ws.Range("E28").Value = "'2026-02-05 17:19:16"
$ws.Range("J28").Value = "'993.8 hPa"
$ws.Range("O28").Value = "'1.9 °C"
$ws.Range("E29").Value = "'2026-02-05 17:19:19"
$ws.Range("H29").Value = "'84%"
$ws.Range("O29").Value = "'8.3 °C"
$ws.Range("E30").Value = "'2026-02-05 17:19:22"
$ws.Range("H30").Value = "'66%"
$ws.Range("K30").Value = "'1.4 MJ/m2"
$ws.Range("O30").Value = "'-2.1 °C"
$ws.Range("E31").Value = "'2026-02-05 17:19:24"
$ws.Range("H31").Value = "'97%"
$ws.Range("I31").Value = "'17.8 mm"
$ws.Range("J31").Value = "'994.7 hPa"
$ws.Range("E32").Value = "'2026-02-05 17:19:27"
$ws.Range("H32").Value = "'83%"
$ws.Range("J32").Value = "'992.0 hPa"
$ws.Range("O32").Value = "'11.6 °C"
$ws.Range("E33").Value = "'2026-02-05 17:19:30"
$ws.Range("O33").Value = "'8.4 °C"
$ws.Range("E34").Value = "'2026-02-05 17:19:33"
$ws.Range("H34").Value = "'99%"
$ws.Range("M34").Value = "'10.4 °C 16:54 TU"
$ws.Range("O34").Value = "'3.0 °C"
$ws.Range("E35").Value = "'2026-02-05 17:19:36"
$ws.Range("I35").Value = "'3.6 mm"
$ws.Range("E36").Value = "'2026-02-05 17:19:38"
$ws.Range("H36").Value = "'86%"
$ws.Range("J36").Value = "'992.5 hPa"
$ws.Range("K36").Value = "'9.7 MJ/m2"
$ws.Range("O36").Value = "'9.9 °C"
